# Auto-generated Excel COM-interop script to apply the country/stat updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 20:37"

# Row 4
$ws.Range("B4").Value = 2172671
$ws.Range("C4").Value = 10443
$ws.Range("D4").Value = 875183
$ws.Range("E4").Value = 1179421
$ws.Range("G4").Value = 209
$ws.Range("H4").Value = 118067

# Row 7
$ws.Range("B7").Value = 342841
$ws.Range("C7").Value = 10058
$ws.Range("D7").Value = 180225
$ws.Range("E7").Value = 152702
$ws.Range("G7").Value = 394
$ws.Range("H7").Value = 9914

# Row 9
$ws.Range("B9").Value = 291189
$ws.Range("C9").Value = 181

# Row 13
$ws.Range("B13").Value = 187927
$ws.Range("C13").Value = 256
$ws.Range("E13").Value = 6452

# Row 14
$ws.Range("A14").Value = "Turquia"
$ws.Range("B14").Value = 179831
$ws.Range("C14").Value = 1592
$ws.Range("D14").Value = 152364
$ws.Range("E14").Value = 22642
$ws.Range("G14").Value = 18
$ws.Range("H14").Value = 4825

# Row 15
$ws.Range("A15").Value = "Chile"
$ws.Range("B15").Value = 179436
$ws.Range("C15").Value = 5143
$ws.Range("D15").Value = 148792
$ws.Range("E15").Value = 27282
$ws.Range("G15").Value = 39
$ws.Range("H15").Value = 3362

# Row 32
$ws.Range("B32").Value = 42636
$ws.Range("C32").Value = 342
$ws.Range("D32").Value = 28129
$ws.Range("E32").Value = 14216
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 291

# Row 49
$ws.Range("B49").Value = 19237
$ws.Range("C49").Value = 182
$ws.Range("D49").Value = 15415
$ws.Range("E49").Value = 3520

# Row 63
$ws.Range("A63").Value = "Azerbaiyan"
$ws.Range("B63").Value = 10324
$ws.Range("C63").Value = 367
$ws.Range("D63").Value = 5739
$ws.Range("E63").Value = 4463
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 122

# Row 64
$ws.Range("A64").Value = "Chequia"
$ws.Range("B64").Value = 10044
$ws.Range("C64").Value = 20
$ws.Range("D64").Value = 7295
$ws.Range("E64").Value = 2419
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 330

# Row 75
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 5263
$ws.Range("C75").Value = 183
$ws.Range("D75").Value = 4019
$ws.Range("E75").Value = 1225
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 19

# Row 76
$ws.Range("A76").Value = "Senegal"
$ws.Range("B76").Value = 5173
$ws.Range("C76").Value = 83
$ws.Range("D76").Value = 3424
$ws.Range("E76").Value = 1685
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 64

# Row 97
$ws.Range("B97").Value = 2310
$ws.Range("C97").Value = 12
$ws.Range("D97").Value = 2058
$ws.Range("E97").Value = 223

# Row 101
$ws.Range("B101").Value = 2065
$ws.Range("C101").Value = 30
$ws.Range("D101").Value = 1540
$ws.Range("E101").Value = 517

# Row 103
$ws.Range("B103").Value = 1905
$ws.Range("C103").Value = 16
$ws.Range("E103").Value = 552

# Row 104
$ws.Range("A104").Value = "Mauritania"
$ws.Range("B104").Value = 1887
$ws.Range("C104").Value = 104
$ws.Range("D104").Value = 360
$ws.Range("E104").Value = 1436
$ws.Range("G104").Value = 4
$ws.Range("H104").Value = 91

# Row 105
$ws.Range("A105").Value = "Mali"
$ws.Range("B105").Value = 1860
$ws.Range("C105").Value = 51
$ws.Range("D105").Value = 1125
$ws.Range("E105").Value = 631
$ws.Range("H105").Value = 104

# Row 106
$ws.Range("A106").Value = "Islandia"
$ws.Range("B106").Value = 1810
$ws.Range("D106").Value = 1796
$ws.Range("E106").Value = 4
$ws.Range("H106").Value = 10

# Row 118
$ws.Range("A118").Value = "Guayana Francesa"
$ws.Range("B118").Value = 1326
$ws.Range("C118").Value = 71
$ws.Range("D118").Value = 552
$ws.Range("E118").Value = 771
$ws.Range("H118").Value = 3

# Row 119
$ws.Range("A119").Value = "Guinea Ecuatorial"
$ws.Range("B119").Value = 1306
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 200
$ws.Range("E119").Value = 1094
$ws.Range("H119").Value = 12

# Row 120
$ws.Range("B120").Value = 1296
$ws.Range("C120").Value = 7
$ws.Range("D120").Value = 673
$ws.Range("E120").Value = 611
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 12

# Row 121
$ws.Range("A121").Value = "Madagascar"
$ws.Range("B121").Value = 1290
$ws.Range("C121").Value = 18
$ws.Range("D121").Value = 384
$ws.Range("E121").Value = 896
$ws.Range("H121").Value = 10

# Row 134
$ws.Range("A134").Value = "Yemen"
$ws.Range("B134").Value = 844
$ws.Range("C134").Value = 116
$ws.Range("D134").Value = 79
$ws.Range("E134").Value = 557
$ws.Range("G134").Value = 44
$ws.Range("H134").Value = 208

# Row 135
$ws.Range("A135").Value = "Cabo Verde"
$ws.Range("B135").Value = 759
$ws.Range("C135").Value = 9
$ws.Range("D135").Value = 301
$ws.Range("E135").Value = 452
$ws.Range("H135").Value = 6

# Row 145
$ws.Range("B145").Value = 555
$ws.Range("C145").Value = 8
$ws.Range("E145").Value = 480

# Row 155
$ws.Range("B155").Value = 387
$ws.Range("C155").Value = 4
$ws.Range("E155").Value = 329

# Row 161
$ws.Range("E161").Value = 195
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 4

# Row 176
$ws.Range("A176").Value = "Eritrea"
$ws.Range("B176").Value = 109
$ws.Range("C176").Value = 13
$ws.Range("D176").Value = 39
$ws.Range("E176").Value = 70
$ws.Range("H176").Value = 0

# Row 177
$ws.Range("A177").Value = "Bahamas"
$ws.Range("B177").Value = 103
$ws.Range("D177").Value = 68
$ws.Range("E177").Value = 24
$ws.Range("H177").Value = 11

# Row 178
$ws.Range("A178").Value = "Aruba"
$ws.Range("B178").Value = 101
$ws.Range("D178").Value = 98
$ws.Range("E178").Value = 0
$ws.Range("H178").Value = 3

# Row 179
$ws.Range("A179").Value = "Monaco"
$ws.Range("B179").Value = 99
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 93
$ws.Range("E179").Value = 2
$ws.Range("H179").Value = 4

# Row 180
$ws.Range("A180").Value = "Barbados"
$ws.Range("B180").Value = 97
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 83
$ws.Range("E180").Value = 7
$ws.Range("H180").Value = 7

# Row 208
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
